# Remove the "is_locked" and "is_enabled" columns (F:G) from the role import
# template's header row. Deleting with a left shift moves the trailing "rem"
# column (H) into F, matching the target layout:
#   A=lbl B=home_url C=menu_ids_lbl D=permit_ids_lbl E=data_permit_ids_lbl F=rem
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlShiftToLeft = -4159
$ws.Range("F1:G1").Delete(-4159)
